$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new data row (time value 0.5 in column A, the string "c" in column B),
# extending the filtered range from A1:B7 to A1:B8.
$ws.Range("A8").Value = 0.5
$ws.Range("B8").Value = "c"

# Grow the existing AutoFilter so its range covers the new row.
$ws.Range("A1:B8").AutoFilter()

# Re-apply the filter on column A (Field 1) with the extra matching value
# (0.500) added alongside the previously selected 0.046 and 0.516, using
# the "filter by list of values" operator (xlFilterValues = 7).
$ws.Range("A1:B8").AutoFilter(1, @("0.046","0.500","0.516"), 7)

# Keep the hidden "_xlnm._FilterDatabase" defined name in sync with the
# grown filter range.
$fdb = $wb.Names.Item(1)
$fdb.RefersTo = "=Munka1!`$A`$1:`$B`$8"

# Move the active cell/selection to where the user ended up after editing.
$ws.Range("C7").Select()
